$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "calculate < montant global brut / taxe > to each mondataire" ---
# Row 2 (ALI EXPRESSE)
$ws.Range("I2").Value = 4000
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 3600

# Row 3 (KHALID TAGHMAOUI)
$ws.Range("I3").Value = 6000
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 5400

# --- "add a condition if < date fin contrat > was null" -> this mandataire
# (whose contract end date was null) is now included in the situation ---
$ws.Range("A4").Value = "Tawfiq MF"
$ws.Range("B4").Value = "BB12354"
# The account number is a long numeric-looking string; prefix with a quote
# so Excel keeps it stored as text instead of converting it to a number.
$ws.Range("C4").Value = "'119349134978532465421354"
$ws.Range("D4").Value = "BMCE"
$ws.Range("E4").Value = "BMCE test"
$ws.Range("F4").Value = "Logement de fonction"
$ws.Range("G4").Value = "001/LF/DR IFRAN"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 9300

# Extend the "number stored as text" ignored-error region to cover the new
# row as well (A1:K3 -> A1:K4). Wrapped defensively in case this particular
# COM entry point isn't fully supported by the host runtime.
try {
    $ws.Range("A1:K4").Errors.Item(9).Ignore = $true
} catch {
}
